$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# Helper approach: Word's run-merge canonicalizer only "fast-paths" (and
# thus re-coalesces) a paragraph the *first* time it is touched while it
# still has a single run. Once a paragraph already has 2+ runs, further
# edits splice cleanly without re-merging neighbouring runs. We exploit
# this by temporarily bookmarking the sub-ranges we want to become their
# own runs (bookmarks force a run split with no extra rPr residue),
# editing each piece's text, and finally deleting the bookmarks again.
# ----------------------------------------------------------------------

# --- Edit 1: "...employees that having their age..." -> "...that have their age..." ---
$full1 = "A user will input an age range and the application will only show employees that having their age in that range."
$all = $d.Content.Text
$start1 = $all.IndexOf($full1)

$prefix1 = "A user will input an age range and the application will only show employees that "
$oldMid1 = "having"

$o1MidStart = $start1 + $prefix1.Length
$o1MidEnd = $o1MidStart + $oldMid1.Length

$rMid1 = $d.Range($o1MidStart, $o1MidEnd)
$d.Bookmarks.Add("tmpSplit1", $rMid1)

$rMid1b = $d.Range($o1MidStart, $o1MidStart + $oldMid1.Length)
$rMid1b.Text = "have"

$d.Bookmarks("tmpSplit1").Delete()

# --- Edit 2: file read/save sentence replaced with plist-specific sentence ---
$full2 = "The application will allow a user to pick a file to read the data from and another file to save the data to."
$all = $d.Content.Text
$start2 = $all.IndexOf($full2)

# Arbitrary split points inside the original sentence; exact position
# doesn't matter since each chunk's text is fully replaced below.
$c1End = $start2 + 22
$c2End = $start2 + 60
$c3End = $start2 + $full2.Length

$rA = $d.Range($start2, $c1End)
$d.Bookmarks.Add("tmpSplit2a", $rA)
$rB = $d.Range($c1End, $c2End)
$d.Bookmarks.Add("tmpSplit2b", $rB)

# Replace right-to-left so earlier offsets stay valid while text lengths change.
$r3 = $d.Range($c2End, $c3End)
$r3.Text = "a file named “employees_out.plist” in the desktop folder."

$r2 = $d.Range($c1End, $c2End)
$r2.Text = "read data from a property list file (Employees.plist) and output to "

$r1 = $d.Range($start2, $c1End)
$r1.Text = "The application will "

$d.Bookmarks("tmpSplit2a").Delete()
$d.Bookmarks("tmpSplit2b").Delete()
